$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCells = @("D2","D3","D5","D8","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D19","D20","D21","D22","D25","D26","D27","D31","D32","D33","D34","D35","D36","D37","D39","D41","D42","D43","D44","D45","D48","D51")
foreach ($addr in $dCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = "34.504.64"
$ws.Range("E2").Value = "  -2.83%  "
$ws.Range("D3").Value = "1.799.73"
$ws.Range("E3").Value = "  -2.26%  "
$ws.Range("E4").Value = "  +0.34%  "
$ws.Range("D5").Value = "229.07"
$ws.Range("E5").Value = "  -1.14%  "
$ws.Range("E6").Value = "  -1.50%  "
$ws.Range("E7").Value = "  +0.37%  "
$ws.Range("D8").Value = "39.15"
$ws.Range("E8").Value = "  -10.93%  "
$ws.Range("D9").Value = "0.319"
$ws.Range("E9").Value = "  +2.85%  "
$ws.Range("D10").Value = "0.0677"
$ws.Range("E10").Value = "  -2.99%  "
$ws.Range("D11").Value = "0.0988"
$ws.Range("E11").Value = "  -2.21%  "
$ws.Range("D12").Value = "2.059.06"
$ws.Range("E12").Value = "  -2.34%  "
$ws.Range("D13").Value = "11.08"
$ws.Range("E13").Value = "  -1.54%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.797.22"
$ws.Range("E14").Value = "  -2.50%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "0.656"
$ws.Range("E15").Value = "  -2.35%  "
$ws.Range("D16").Value = "4.54"
$ws.Range("E16").Value = "  -3.61%  "
$ws.Range("D17").Value = "34.328.65"
$ws.Range("E17").Value = "  -3.30%  "
$ws.Range("D18").Value = "68.96"
$ws.Range("E18").Value = "  -2.11%  "
$ws.Range("D19").Value = "0.0₃0777"
$ws.Range("E19").Value = "  -2.89%  "
$ws.Range("D20").Value = "239.40"
$ws.Range("E20").Value = "  -2.02%  "
$ws.Range("D21").Value = "11.76"
$ws.Range("E21").Value = "  -2.43%  "
$ws.Range("D22").Value = "4.67"
$ws.Range("E22").Value = "  +0.80%  "
$ws.Range("E23").Value = "  +0.42%  "
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").Value = "173.08"
$ws.Range("E25").Value = "  +1.08%  "
$ws.Range("D26").Value = "7.68"
$ws.Range("E26").Value = "  -3.91%  "
$ws.Range("D27").Value = "17.16"
$ws.Range("E27").Value = "  -3.65%  "
$ws.Range("E28").Value = "  -0.44%  "
$ws.Range("E29").Value = "  -4.86%  "
$ws.Range("E30").Value = "  +0.39%  "
$ws.Range("D31").Value = "4.00"
$ws.Range("E31").Value = "  +1.54%  "
$ws.Range("D32").Value = "0.0540"
$ws.Range("E32").Value = "  -1.97%  "
$ws.Range("D33").Value = "3.89"
$ws.Range("E33").Value = "  -4.92%  "
$ws.Range("D34").Value = "1.21"
$ws.Range("E34").Value = "  +7.37%  "
$ws.Range("D35").Value = "1.79"
$ws.Range("E35").Value = "  -2.73%  "
$ws.Range("D36").Value = "0.690"
$ws.Range("E36").Value = "  -0.25%  "
$ws.Range("D37").Value = "90.57"
$ws.Range("E37").Value = "  -4.44%  "
$ws.Range("E38").Value = "  +4.63%  "
$ws.Range("D39").Value = "1.320.61"
$ws.Range("E39").Value = "  -1.99%  "
$ws.Range("E40").Value = "  -2.67%  "
$ws.Range("D41").Value = "0.955"
$ws.Range("E41").Value = "  -5.58%  "
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").Value = "14.22"
$ws.Range("E42").Value = "  -7.50%  "
$ws.Range("B43").Value = "HuobiToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D43").Value = "2.40"
$ws.Range("E43").Value = "  -2.32%  "
$ws.Range("D44").Value = "2.20"
$ws.Range("E44").Value = "  -9.75%  "
$ws.Range("D45").Value = "2.70"
$ws.Range("E45").Value = "  -3.70%  "
$ws.Range("E46").Value = "  -1.94%  "
$ws.Range("E47").Value = "  -1.34%  "
$ws.Range("D48").Value = "1.981.18"
$ws.Range("E48").Value = "  -1.55%  "
$ws.Range("E49").Value = "  +0.35%  "
$ws.Range("E50").Value = "  +3.63%  "
$ws.Range("D51").Value = "97.50"
$ws.Range("E51").Value = "  -4.78%  "
